# Textbox response formatting fix
# Rename sheets and update stimulus/response CSV filename values (timestamps refreshed).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511686477695434"
$ws1.Range("B2").Value = "go_stims-16511686477375429.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686477535434.csv"
$ws1.Range("B4").Value = "go_stims-16511686477545457.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686477685437.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511686501867845"
$ws2.Range("B2").Value = "OB-16511686492647562.csv"
$ws2.Range("B3").Value = "OB-1651168649015754.csv"
$ws2.Range("B4").Value = "ZB-match_2-16511686477935436.csv"
$ws2.Range("B5").Value = "TB-16511686501657517.csv"
$ws2.Range("B6").Value = "OB-16511686485897553.csv"
$ws2.Range("B7").Value = "TB-16511686494227543.csv"
$ws2.Range("B8").Value = "TB-16511686494937856.csv"
$ws2.Range("B9").Value = "ZB-match_0-16511686478855445.csv"
$ws2.Range("B10").Value = "ZB-match_5-16511686481607547.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16511686501877522"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511686502347858"
$ws4.Range("B2").Value = "MM_stims-1651168650202785.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686501897526.csv"
$ws4.Range("B4").Value = "MM_stims-16511686502187657.csv"
$ws4.Range("B5").Value = "ZM_stims-1651168650202785.csv"
$ws4.Range("B6").Value = "MM_stims-16511686502347858.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686502197542.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511686502976882"
$ws5.Range("B2").Value = "SAT_stims-16511686502507863.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686502664113.csv"
$ws5.Range("B4").Value = "SAT_stims-16511686502377517.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511686502824118.csv"
